$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DATE_TYPE_CODE (keep stored as text, like the original "004")
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "002"
$ws.Range("J2").Style = "Normal"

# NOTICE_DATE / REPORT_DATE (stored as text strings, not real dates)
$ws.Range("M2").Value = "2020-12-23 00:00:00"
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# TOTAL_ASSETS
$ws.Range("O2").Value = 786079324.92
# FIXED_ASSET
$ws.Range("P2").Value = 15568801.25
# MONETARYFUNDS
$ws.Range("Q2").Value = 135638478.24
# MONETARYFUNDS_RATIO
$ws.Range("R2").Value = 16.0342927059
# ACCOUNTS_RECE
$ws.Range("S2").Value = 527638005.02
# ACCOUNTS_RECE_RATIO
$ws.Range("T2").Value = 45.0499947807

# TOTAL_LIABILITIES
$ws.Range("W2").Value = 325856225.71
# ACCOUNTS_PAYABLE
$ws.Range("X2").Value = 13303382.44
# ACCOUNTS_PAYABLE_RATIO
$ws.Range("Y2").Value = 176.9790674346

# TOTAL_EQUITY
$ws.Range("AB2").Value = 460223099.21
# TOTAL_EQUITY_RATIO
$ws.Range("AC2").Value = 27.7101796972
# TOTAL_ASSETS_RATIO
$ws.Range("AD2").Value = 42.1393293087
# TOTAL_LIAB_RATIO
$ws.Range("AE2").Value = 69.1273810531

# CURRENT_RATIO
$ws.Range("AF2").Value = 234.2191464938
# DEBT_ASSET_RATIO
$ws.Range("AG2").Value = 41.4533515104
